$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column H (the "pt_max" column, a constant 50 in every data row).
# This shifts everything after it left by one, and Excel auto-adjusts
# formula references and the used-range dimension.
$ws.Range("H:H").Delete()

# Header row 1 becomes bold.
$ws.Range("A1:Q1").Font.Bold = $true

# Move the selection cursor, matching the saved cursor position in the file.
$ws.Range("O12").Select() | Out-Null
